# Insert 3 new data rows before the current row 33, shifting the existing
# rows 33-111 down to 36-114, then populate the 3 newly inserted rows with
# their data (rows 33, 34 and 35 in the final sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows starting at row 33.
$ws.Rows("33:35").Insert()

# Common (constant) values shared by every Damasco data row in this sheet.
$marketId = 10
$market = "Vega Modelo de Temuco"
$region = "La Araucanía"
$codreg = 9
$tipo = "Fruta"
$prodId = 100103
$producto = "Frutos de hueso (carozo)"
$catId = 100103003
$categoria = "Damasco"

# ----- Row 33 -----
$ws.Range("A33").Value2 = $marketId
$ws.Range("B33").Value2 = $market
$ws.Range("C33").Value2 = $region
$ws.Range("D33").Value2 = 45281
$ws.Range("D33").NumberFormat = $ws.Range("D36").NumberFormat
$ws.Range("E33").Value2 = $codreg
$ws.Range("F33").Value2 = $tipo
$ws.Range("G33").Value2 = $prodId
$ws.Range("H33").Value2 = $producto
$ws.Range("I33").Value2 = $catId
$ws.Range("J33").Value2 = $categoria
$ws.Range("K33").Value2 = "Castle Brite"
$ws.Range("L33").Value2 = "Especial"
$ws.Range("M33").Value2 = 80
$ws.Range("N33").Value2 = 25000
$ws.Range("O33").Value2 = 25000
$ws.Range("P33").Value2 = 25000
$ws.Range("Q33").Value2 = "$/bandeja 18 kilos"
$ws.Range("R33").Value2 = "Región de O'Higgins"
$ws.Range("S33").Value2 = 1389
$ws.Range("T33").Value2 = 18

# ----- Row 34 -----
$ws.Range("A34").Value2 = $marketId
$ws.Range("B34").Value2 = $market
$ws.Range("C34").Value2 = $region
$ws.Range("D34").Value2 = 45281
$ws.Range("D34").NumberFormat = $ws.Range("D36").NumberFormat
$ws.Range("E34").Value2 = $codreg
$ws.Range("F34").Value2 = $tipo
$ws.Range("G34").Value2 = $prodId
$ws.Range("H34").Value2 = $producto
$ws.Range("I34").Value2 = $catId
$ws.Range("J34").Value2 = $categoria
$ws.Range("K34").Value2 = "Castle Brite"
$ws.Range("L34").Value2 = "Primera"
$ws.Range("M34").Value2 = 550
$ws.Range("N34").Value2 = 17000
$ws.Range("O34").Value2 = 18000
$ws.Range("P34").Value2 = 17545
$ws.Range("Q34").Value2 = "$/caja 15 kilos"
$ws.Range("R34").Value2 = "Región de O'Higgins"
$ws.Range("S34").Value2 = 1170
$ws.Range("T34").Value2 = 15

# ----- Row 35 -----
$ws.Range("A35").Value2 = $marketId
$ws.Range("B35").Value2 = $market
$ws.Range("C35").Value2 = $region
$ws.Range("D35").Value2 = 45281
$ws.Range("D35").NumberFormat = $ws.Range("D36").NumberFormat
$ws.Range("E35").Value2 = $codreg
$ws.Range("F35").Value2 = $tipo
$ws.Range("G35").Value2 = $prodId
$ws.Range("H35").Value2 = $producto
$ws.Range("I35").Value2 = $catId
$ws.Range("J35").Value2 = $categoria
$ws.Range("K35").Value2 = "Castle Brite"
$ws.Range("L35").Value2 = "Segunda"
$ws.Range("M35").Value2 = 100
$ws.Range("N35").Value2 = 15000
$ws.Range("O35").Value2 = 15000
$ws.Range("P35").Value2 = 15000
$ws.Range("Q35").Value2 = "$/caja 15 kilos"
$ws.Range("R35").Value2 = "Región de O'Higgins"
$ws.Range("S35").Value2 = 1000
$ws.Range("T35").Value2 = 15
